# Natmi following Dr Hou advice
# Adds the "M2" macrophage cluster to the sending-cluster side of the
# Igf2 -> Igf2r ligand-receptor communication table (previously only a
# target/receiving cluster) and refreshes all computed NATMI metrics
# for the resulting 4x4 sending x target cluster grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($sheet, $rowNum, $values)
    $arr = New-Object "object[,]" 1,20
    for ($i = 0; $i -lt 20; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rng = $sheet.Range("A" + $rowNum + ":T" + $rowNum)
    $rng.Value = $arr
}

Set-RowValues $ws 2 @("ECs", "Igf2", "Igf2r", "ECs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 18.86633166666667, 56.598995, 0.1413973975846522, 0.1413973975846523, 161.1875106296789, 1450.68759566711, 0.003605156542418446, 0.003605156542418447)
Set-RowValues $ws 3 @("ECs", "Igf2", "Igf2r", "FAPs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 81.49602766666666, 244.488083, 0.6107878536829223, 0.6107878536829223, 696.2742974039081, 6266.468676635173, 0.01557302937924594, 0.01557302937924594)
Set-RowValues $ws 4 @("ECs", "Igf2", "Igf2r", "M2", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 10.11432133333333, 30.342964, 0.07580375136704795, 0.07580375136704795, 86.41331585986578, 777.7198427387921, 0.001932739886652888, 0.001932739886652888)
Set-RowValues $ws 5 @("ECs", "Igf2", "Igf2r", "sCs", 3, 1, 8.543659333333332, 25.630978, 0.02549662585027493, 0.02549662585027493, 3, 1, 22.95103433333334, 68.853103, 0.1720109973653774, 0.1720109973653775, 196.0858186916371, 1764.772368224734, 0.004385700041957655, 0.004385700041957656)
Set-RowValues $ws 6 @("FAPs", "Igf2", "Igf2r", "ECs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 18.86633166666667, 56.598995, 0.1413973975846522, 0.1413973975846523, 5137.828797925695, 46240.45918133125, 0.1149138480537908, 0.1149138480537909)
Set-RowValues $ws 7 @("FAPs", "Igf2", "Igf2r", "FAPs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 81.49602766666666, 244.488083, 0.6107878536829223, 0.6107878536829223, 22193.64343107236, 199742.7908796512, 0.4963880793435397, 0.4963880793435396)
Set-RowValues $ws 8 @("FAPs", "Igf2", "Igf2r", "M2", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 10.11432133333333, 30.342964, 0.07580375136704795, 0.07580375136704795, 2754.412057203889, 24789.708514835, 0.06160580686278345, 0.06160580686278345)
Set-RowValues $ws 9 @("FAPs", "Igf2", "Igf2r", "sCs", 3, 1, 272.3279166666667, 816.98375, 0.8127012944845317, 0.8127012944845315, 3, 1, 22.95103433333334, 68.853103, 0.1720109973653774, 0.1720109973653775, 6250.207365341807, 56251.86628807625, 0.1397935602244176, 0.1397935602244176)
Set-RowValues $ws 10 @("M2", "Igf2", "Igf2r", "ECs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 18.86633166666667, 56.598995, 0.1413973975846522, 0.1413973975846523, 50.11517135223055, 451.036542170075, 0.00112088732662423, 0.001120887326624231)
Set-RowValues $ws 11 @("M2", "Igf2", "Igf2r", "FAPs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 81.49602766666666, 244.488083, 0.6107878536829223, 0.6107878536829223, 216.4802073450839, 1948.321866105755, 0.004841845579508132, 0.004841845579508132)
Set-RowValues $ws 12 @("M2", "Igf2", "Igf2r", "M2", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 10.11432133333333, 30.342964, 0.07580375136704795, 0.07580375136704795, 26.86695833017111, 241.80262497154, 0.0006009125038318307, 0.0006009125038318307)
Set-RowValues $ws 13 @("M2", "Igf2", "Igf2r", "sCs", 2, 0.6666666666666666, 2.656328333333333, 7.968985, 0.007927213271044639, 0.007927213271044639, 3, 1, 22.95103433333334, 68.853103, 0.1720109973653774, 0.1720109973653775, 60.96548277893945, 548.689345010455, 0.001363567861080445, 0.001363567861080445)
Set-RowValues $ws 14 @("sCs", "Igf2", "Igf2r", "ECs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 18.86633166666667, 56.598995, 0.1413973975846522, 0.1413973975846523, 972.7838816083912, 8755.05493447552, 0.02175750566181871, 0.02175750566181871)
Set-RowValues $ws 15 @("sCs", "Igf2", "Igf2r", "FAPs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 81.49602766666666, 244.488083, 0.6107878536829223, 0.6107878536829223, 4202.08992028453, 37818.80928256077, 0.09398489938062865, 0.09398489938062862)
Set-RowValues $ws 16 @("sCs", "Igf2", "Igf2r", "M2", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 10.11432133333333, 30.342964, 0.07580375136704795, 0.07580375136704795, 521.5136116714383, 4693.622505042945, 0.01166429211377979, 0.01166429211377978)
Set-RowValues $ws 17 @("sCs", "Igf2", "Igf2r", "sCs", 3, 1, 51.56189866666667, 154.685696, 0.1538748663941489, 0.1538748663941489, 3, 1, 22.95103433333334, 68.853103, 0.1720109973653774, 0.1720109973653775, 1183.398906590521, 10650.59015931469, 0.02646816923792175, 0.02646816923792175)
